# updated legacy GSC export data
# Appends 4 new daily rows (2025-11-12 .. 2025-11-15) to the "Chart" sheet
# and corrects the "Page with redirect" Pages count on the
# "Critical issues" sheet from 22 to 24.

$wb = $excel.ActiveWorkbook

# ---- Chart sheet: append rows 86-89 --------------------------------------
$chart = $wb.Worksheets.Item("Chart")

$newRows = @(
    @{ Row = 86; Date = "2025-11-12"; NotIndexed = 38; Indexed = 97; Impressions = 16 },
    @{ Row = 87; Date = "2025-11-13"; NotIndexed = 38; Indexed = 97; Impressions = 11 },
    @{ Row = 88; Date = "2025-11-14"; NotIndexed = 38; Indexed = 97; Impressions = 10 },
    @{ Row = 89; Date = "2025-11-15"; NotIndexed = 38; Indexed = 97; Impressions = 8 }
)

foreach ($r in $newRows) {
    # Leading apostrophe forces the date-like text to be stored as a plain
    # string (matching every other date cell in the column) instead of
    # being auto-converted to a date serial number.
    $chart.Cells.Item($r.Row, 1).Value = "'" + $r.Date
    $chart.Cells.Item($r.Row, 2).Value = $r.NotIndexed
    $chart.Cells.Item($r.Row, 3).Value = $r.Indexed
    $chart.Cells.Item($r.Row, 4).Value = $r.Impressions
}

# ---- Critical issues sheet: "Page with redirect" Pages 22 -> 24 ----------
$critical = $wb.Worksheets.Item("Critical issues")
$critical.Range("D3").Value = 24

Write-Host "Applied GSC export update: added rows 86-89 to Chart, updated Critical issues D3 to 24"
